$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "daily events" rows appended below the header (row 1). Columns D/E/F
# (vehicleId / Unidad / driverId) look numeric but must stay text, so each
# is briefly forced to the "@" (Text) number format before the value is
# typed in, then reset back to the default "Normal" style so the saved
# cell carries no explicit style index (matching the source data which was
# written as plain inline strings).

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 1).Value = "281474991395097-1750090341759"
$ws.Cells.Item(2, 2).Value = "Harsh Brake"
$ws.Cells.Item(2, 3).Value = "2025-06-16T10:12:21.759"
Set-TextCell 2 4 "281474991395097"
Set-TextCell 2 5 "125"
Set-TextCell 2 6 "51834055"
$ws.Cells.Item(2, 7).Value = "DAVID SERRANO"
$ws.Cells.Item(2, 8).Value = 20.597474099
$ws.Cells.Item(2, 9).Value = -103.43804784
$ws.Cells.Item(2, 10).Value = 0.8824904561042786
$ws.Cells.Item(2, 11).Value = "No video URL"
$ws.Cells.Item(2, 12).Value = "No video URL"

# Row 3
$ws.Cells.Item(3, 1).Value = "281474991395097-1750090239896"
$ws.Cells.Item(3, 2).Value = "Harsh Brake"
$ws.Cells.Item(3, 3).Value = "2025-06-16T10:10:39.896"
Set-TextCell 3 4 "281474991395097"
Set-TextCell 3 5 "125"
Set-TextCell 3 6 "51834055"
$ws.Cells.Item(3, 7).Value = "DAVID SERRANO"
$ws.Cells.Item(3, 8).Value = 20.60085944
$ws.Cells.Item(3, 9).Value = -103.43621207
$ws.Cells.Item(3, 10).Value = 0.7449945211410522
$ws.Cells.Item(3, 11).Value = "No video URL"
$ws.Cells.Item(3, 12).Value = "No video URL"

# Row 4
$ws.Cells.Item(4, 1).Value = "281474991205821-1750089224772"
$ws.Cells.Item(4, 2).Value = "No Seat Belt"
$ws.Cells.Item(4, 3).Value = "2025-06-16T09:53:44.772"
Set-TextCell 4 4 "281474991205821"
Set-TextCell 4 5 "148"
Set-TextCell 4 6 "51834015"
$ws.Cells.Item(4, 7).Value = "LUIS IBARRA"
$ws.Cells.Item(4, 8).Value = 20.63196654
$ws.Cells.Item(4, 9).Value = -103.436270769
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1750089222272/NAzsWI3hGH-camera-video-segment-driver-1750089224772.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBDGJV2O2%2F20250617%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250617T170515Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDLaChzm8LA2YajoH1jmoheQifOTueA9A2DH1awyWsPegIgVbOJNH7YbydftCPqFkRE7aS263fdvF4RyKWitA0v0vwq3QMIdRAEGgw3ODEyMDQ5NDIyNDQiDLewVzCH%2BNkMnCxsQyq6A5BnEj4YLA%2BjqTFC%2BuexdMeHT4F%2Fkod694yR9hG5v8tWy8XQYGDSjZTyCvb40ZCPCx2iSPmbYcppkaxuXEidnGHAdVHfPAX%2BcZ%2BGEeTeMwQGGNP9B0gtcVhMpF0C8q%2FaUp%2BGqLKTY06pA22%2B8pZ%2Bm72TkPHCRocWHnFtNj%2F8PctKY4gDSZLdsC%2B2z75qGzvK7t7qbrFXV6JI0PDv0ICb2lXTbbPgE7ghMwweWLjYz6n1AWmdpQD7MKSqP9n5w4O3wKRYwkFW4n4yBu7sr6m5DreDPFwFJ2otul%2Bigw1Y%2BYcPUL0KtZagcyk8cVO1VF1MDIndb70O0YqfT6F9qkUXwP75hp0HFfeVbFD%2B8jU86QSnwfpbnMKEudcg0CD%2BMfwyKSkkrc6xD3HyahvYoSDsqkesHCDf8UD0Jm1%2BnqqqgMdPkmJVV67nLgD6j3CT5WaKAuNtdqVKnF%2BWDQsN0CRjzbL84QK%2FkHiOlLjkPmdvKh0sYimQqtbEFk%2FRq8o9pkfsqxZTjxyz9nZ2dr3efslVRtcDegEHgHZGOkITTVVsIjl%2FxcpcrB%2FTP5ofG81feb7YPBrAXEVLHdIHD0Ewk6rFwgY6pQEGLXZe%2BBq8zneBvk8ROnDYskIOSd24ladlrDVs5SOl0wGZiK8HZb2UpHr4BkaHwq4sPEdj6pNvtlcHy%2BbL6y17k%2Bve7HM%2ByVLNXREF8a1ZLlLGrYXghWiWhUyrt5ek4xLbGA640Z4iWVj3SyCKoy8wKHct943FSnaaLj5pBRxQse%2FNVfkumtuf7sfsqLG5Ipf3xmi9eeMrjrB%2Fx4YLTKIG3O4x6zM%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2018%20Jun%202025%2001%3A05%3A15%20GMT&X-Amz-Signature=f6d62db13c039551a6b2ebbb02b73094f3edf973ba4e4674878ce91d99376eeb"
$ws.Cells.Item(4, 12).Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1750089222272/DaHDiAhqBq-camera-video-segment-1750089224772.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBDGJV2O2%2F20250617%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250617T170515Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDLaChzm8LA2YajoH1jmoheQifOTueA9A2DH1awyWsPegIgVbOJNH7YbydftCPqFkRE7aS263fdvF4RyKWitA0v0vwq3QMIdRAEGgw3ODEyMDQ5NDIyNDQiDLewVzCH%2BNkMnCxsQyq6A5BnEj4YLA%2BjqTFC%2BuexdMeHT4F%2Fkod694yR9hG5v8tWy8XQYGDSjZTyCvb40ZCPCx2iSPmbYcppkaxuXEidnGHAdVHfPAX%2BcZ%2BGEeTeMwQGGNP9B0gtcVhMpF0C8q%2FaUp%2BGqLKTY06pA22%2B8pZ%2Bm72TkPHCRocWHnFtNj%2F8PctKY4gDSZLdsC%2B2z75qGzvK7t7qbrFXV6JI0PDv0ICb2lXTbbPgE7ghMwweWLjYz6n1AWmdpQD7MKSqP9n5w4O3wKRYwkFW4n4yBu7sr6m5DreDPFwFJ2otul%2Bigw1Y%2BYcPUL0KtZagcyk8cVO1VF1MDIndb70O0YqfT6F9qkUXwP75hp0HFfeVbFD%2B8jU86QSnwfpbnMKEudcg0CD%2BMfwyKSkkrc6xD3HyahvYoSDsqkesHCDf8UD0Jm1%2BnqqqgMdPkmJVV67nLgD6j3CT5WaKAuNtdqVKnF%2BWDQsN0CRjzbL84QK%2FkHiOlLjkPmdvKh0sYimQqtbEFk%2FRq8o9pkfsqxZTjxyz9nZ2dr3efslVRtcDegEHgHZGOkITTVVsIjl%2FxcpcrB%2FTP5ofG81feb7YPBrAXEVLHdIHD0Ewk6rFwgY6pQEGLXZe%2BBq8zneBvk8ROnDYskIOSd24ladlrDVs5SOl0wGZiK8HZb2UpHr4BkaHwq4sPEdj6pNvtlcHy%2BbL6y17k%2Bve7HM%2ByVLNXREF8a1ZLlLGrYXghWiWhUyrt5ek4xLbGA640Z4iWVj3SyCKoy8wKHct943FSnaaLj5pBRxQse%2FNVfkumtuf7sfsqLG5Ipf3xmi9eeMrjrB%2Fx4YLTKIG3O4x6zM%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2018%20Jun%202025%2001%3A05%3A15%20GMT&X-Amz-Signature=0afd0a119805c8e7773a3013df62097feb6e0e5d5857f45f40a142be8cf8b9b7"
